$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 14
$ws.Range("A3").Value = 63.60000000000036
$ws.Range("A4").Value = 1
$ws.Range("A7").Value = 43.20000000000073
$ws.Range("A8").Value = 8.800000000001091
$ws.Range("A9").Value = 11
